# Generate Report for Handoff
# Update the localization-status report with the latest handoff run:
#  - Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#    newly generated handoff batch (rows 8-13).
#  - zh-cn / de-de sheets: mark the corresponding rows as handoff type
#    "ht" in the Priority column and refresh "Latest Handoff Datetime".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" (column G) for rows 8-13.
$wsOverview.Range("G8:G13").Value = "2016-08-12 02:43:46"

# zh-cn: "Priority" (column E) and "Latest Handoff Datetime" (column H)
# for rows 8-13.
$wsZhCn.Range("E8:E13").Value = "ht"
$wsZhCn.Range("H8:H13").Value = "2016-08-12 02:43:41"

# de-de: "Priority" (column E) for rows 8-13.
$wsDeDe.Range("E8:E13").Value = "ht"
